$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F holds the "2-2" attendance marks for each student row (4-9).
# Reset those marks from 1 to 0.
$ws.Range("F4:F9").Value = 0
